# Add "Chenuis Falls Trail" as a new row in the Hike Difficulties table.
# The row is inserted as the new row 8 (between "Carbon River to Ipsut
# Falls" and "Crystal Lakes and Sourdough Gap"), shifting every
# subsequent row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# There is a single Excel Table ("Table1") covering A1:D64 - grab a
# reference to it before we start moving rows around.
$lo = $ws.ListObjects.Item(1)

# Shift rows 8..64 down to 9..65, opening up a blank row 8.
$ws.Rows("8:8").Insert()

# Fill in the data for the newly-added hike.
$ws.Range("A8").Value = "Chenuis Falls Trail"
$ws.Range("B8").Value = 10
$ws.Range("C8").Value = 1270
$ws.Range("D8").Value = "moderate"

# Grow the table/autofilter range so it covers the new row as well.
$lo.Resize($ws.Range("A1:D65"))

# Match the saved selection from the edited workbook.
$ws.Range("D9").Select()
